# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2-24) a new leading value is inserted into column B,
# pushing the existing values in B:J one column to the right (into C:K).
# The value that previously sat in column K is discarded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing values one column to the right: B2:J24 -> C2:K24
# (this naturally overwrites/discards whatever used to be in column K)
$src = $ws.Range("B2:J24")
$dst = $ws.Range("C2:K24")
$dst.Value = $src.Value()

# New values to insert into column B for each row
$newValues = @{
    2  = 2.48932270964054
    3  = 9.643547872076862
    4  = -9.584088888243137
    5  = -0.573770965293057
    6  = 1.563148290176452
    7  = -1.985496228563019
    8  = -1.838568686009481
    9  = 0.319385421520574
    10 = -0.2177157015159319
    11 = -0.1395947820665385
    12 = -0.3119065001142551
    13 = 0.7021231295320197
    14 = 1.514070997382048
    15 = 0.2163102553365951
    16 = 0.3684555432821496
    17 = 0.661541622456546
    18 = -0.07992401592518952
    19 = 0.1551026493581833
    20 = -0.08373363042288225
    21 = 0.1925427069667326
    22 = -0.4379379024501944
    23 = 0.2324016585002178
    24 = -0.09587373626955231
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
